$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Set up template cell far away with desired font/format, then copy style to destination.
$tmpl = $ws.Range("Z1")
$tmpl.Font.Name = "Arial"
$tmpl.Font.Size = 10
$tmpl.Font.ColorIndex = 8
$tmpl.NumberFormat = "#,##0.0;-#,##0.0;-"
$tmpl.Borders.LineStyle = 1
$tmpl.HorizontalAlignment = -4152
$tmpl.Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
